$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B / C (coin name / link) swaps ---
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"

# --- Column D (price) updates ---
# Cells whose new value would be auto-parsed as a genuine number by Excel
# get their NumberFormat forced to Text ("@") first, so the literal string
# (with its exact digits/trailing zeros) is preserved, matching the source
# workbook where every Price cell is stored as inline text.
$ws.Range("D2").Value = "67.680.65"
$ws.Range("D3").Value = "3.772.51"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.02"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.16"
$ws.Range("D7").Value = "3.768.10"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.90"
$ws.Range("D15").Value = "4.408.02"
$ws.Range("D16").Value = "3.791.52"
$ws.Range("D17").Value = "67.685.38"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.11"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "455.13"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.41"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.688"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.78"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.80"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.07"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.87"
$ws.Range("D30").Value = "3.926.01"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.20"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.57"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.72"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0988"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.145"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.76"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.974"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.13"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.46"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.11"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.293"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.60"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.34"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.82"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "381.57"

# --- Column E (1h volume/percent change) updates ---
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  -3.22%  "
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("E10").Value = "  -3.48%  "
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("E12").Value = "  +4.76%  "
$ws.Range("E13").Value = "  -4.28%  "
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("E22").Value = "  -4.85%  "
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("E25").Value = "  -7.27%  "
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("E31").Value = "  -2.89%  "
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("E33").Value = "  -7.08%  "
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  -1.75%  "
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("E38").Value = "  +5.00%  "
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -8.75%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("E51").Value = "  -2.55%  "
